# إضافة حدث جديد في Card21 by admin at 2025-12-08 08:43:31
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card21")

# Row 19 currently has empty inline-string placeholders in columns B:K.
# Fill them with the literal text "nan" to match the rest of the sheet's
# convention for missing values.
$ws.Range("B19:K19").Value = "nan"

# Append a brand-new service event as row 20 (columns B:K stay blank,
# matching how a brand-new row is first written before being backfilled).
# Copy A19 (already the text "21") down into A20 so the card number lands
# as text rather than being auto-coerced to a number.
$ws.Range("A19").Copy($ws.Range("A20"))
$ws.Range("L20").Value = "14\8\2025"
$ws.Range("M20").Value = "9754 h   770 t"
$ws.Range("N20").Value = "تم تغيير زيت الجيربوكس"
$ws.Range("O20").Value = "تيم العمل"
